# Elimna EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "database" feeding the account-statement table (rows 16-29, columns
# B:G) is re-sorted: instead of alternating between the two workers for
# each period, the rows are now grouped by worker (CLAUDIA first, then
# PATRICIA) and, within each worker, ordered by period descending
# (1910 -> 1904). The underlying (worker, period, amount) data is the
# same 14-row set, just reshuffled into the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row order for B16:G29 -> Tipo, NumDoc, Nombre, Periodo, ValorMora, Salario
$rows = @(
    @("CC","45592819","CLAUDIA CECILIA GARCIA MONCRIEFF","1910",23187,828116),
    @("CC","45592819","CLAUDIA CECILIA GARCIA MONCRIEFF","1909",33125,828116),
    @("CC","45592819","CLAUDIA CECILIA GARCIA MONCRIEFF","1908",33125,828116),
    @("CC","45592819","CLAUDIA CECILIA GARCIA MONCRIEFF","1907",33125,828116),
    @("CC","45592819","CLAUDIA CECILIA GARCIA MONCRIEFF","1906",33125,828116),
    @("CC","45592819","CLAUDIA CECILIA GARCIA MONCRIEFF","1905",33125,828116),
    @("CC","45592819","CLAUDIA CECILIA GARCIA MONCRIEFF","1904",33125,828116),
    @("CC","45592816","PATRICIA MARGARITA GARCIA MONCRIEFF","1910",21874,828116),
    @("CC","45592816","PATRICIA MARGARITA GARCIA MONCRIEFF","1909",31249,828116),
    @("CC","45592816","PATRICIA MARGARITA GARCIA MONCRIEFF","1908",31249,828116),
    @("CC","45592816","PATRICIA MARGARITA GARCIA MONCRIEFF","1907",31249,828116),
    @("CC","45592816","PATRICIA MARGARITA GARCIA MONCRIEFF","1906",31249,828116),
    @("CC","45592816","PATRICIA MARGARITA GARCIA MONCRIEFF","1905",33125,828116),
    @("CC","45592816","PATRICIA MARGARITA GARCIA MONCRIEFF","1904",33125,828116)
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Range("B$r").Value = $data[0]
    $ws.Range("C$r").Value = $data[1]
    $ws.Range("D$r").Value = $data[2]
    $ws.Range("E$r").Value = $data[3]
    $ws.Range("F$r").Value = $data[4]
    $ws.Range("G$r").Value = $data[5]
}
